# Auto-generated edit script applying cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng, [string]$val) {
    # Force text storage so numeric-looking strings (e.g. "6.02", "0.614")
    # are not reinterpreted by Excel as numbers, while avoiding a lingering
    # non-default cell style by resetting back to the Normal style afterward.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "37.867.85"
$ws.Range("E2").Value = "  -0.37%  "
Set-TextCell $ws.Range("D3") "2.030.77"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextCell $ws.Range("D5") "227.28"
$ws.Range("E5").Value = "  -1.15%  "
Set-TextCell $ws.Range("D6") "0.614"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +2.34%  "
Set-TextCell $ws.Range("D9") "0.385"
$ws.Range("E9").Value = "  -0.53%  "
Set-TextCell $ws.Range("D10") "0.0811"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +0.30%  "
Set-TextCell $ws.Range("D12") "14.64"
$ws.Range("E12").Value = "  +0.17%  "
Set-TextCell $ws.Range("D13") "2.331.72"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  +2.43%  "
Set-TextCell $ws.Range("D15") "0.765"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  -1.79%  "
Set-TextCell $ws.Range("D17") "2.028.03"
$ws.Range("E17").Value = "  -1.40%  "
Set-TextCell $ws.Range("D18") "37.760.64"
$ws.Range("E18").Value = "  -0.43%  "
Set-TextCell $ws.Range("D19") "6.02"
$ws.Range("E19").Value = "  -1.85%  "
Set-TextCell $ws.Range("D20") "69.95"
$ws.Range("E20").Value = "  +0.26%  "
Set-TextCell $ws.Range("D21") "0.0₃0823"
$ws.Range("E21").Value = "  -0.88%  "
Set-TextCell $ws.Range("D22") "225.01"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextCell $ws.Range("D24") "2.40"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  -0.11%  "
Set-TextCell $ws.Range("D27") "165.04"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -2.98%  "
Set-TextCell $ws.Range("D29") "18.94"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  -4.80%  "
Set-TextCell $ws.Range("D31") "0.118"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Range("D34") "0.0602"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D35") "4.49"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E37").Value = "  -3.69%  "
Set-TextCell $ws.Range("D38") "3.23"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("E39").Value = "  +0.05%  "
Set-TextCell $ws.Range("D40") "1.524.26"
$ws.Range("E40").Value = "  +2.58%  "
Set-TextCell $ws.Range("D41") "0.0218"
$ws.Range("E41").Value = "  -0.01%  "
Set-TextCell $ws.Range("D42") "96.62"
$ws.Range("E42").Value = "  -1.49%  "
Set-TextCell $ws.Range("D43") "16.80"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("E44").Value = "  -0.56%  "
Set-TextCell $ws.Range("D45") "0.0918"
$ws.Range("E45").Value = "  -1.80%  "
Set-TextCell $ws.Range("D46") "4.11"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -0.28%  "
Set-TextCell $ws.Range("D51") "2.219.58"
$ws.Range("E51").Value = "  -1.26%  "
